$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FedEx ShipmentTracking numbers for rows 2-25 (column P), dated 1 Apr 2022.
$newTracking = @(
    "320018208097",
    "320018208101",
    "320018208134",
    "320018208156",
    "320018208190",
    "320018208215",
    "320018208248",
    "320018208260",
    "320018208292",
    "320018208318",
    "320018208351",
    "320018208373",
    "320018208400",
    "320018208421",
    "320018208454",
    "320018208476",
    "320018208513",
    "320018208535",
    "320018208568",
    "320018208580",
    "320018208616",
    "320018208627",
    "320018208638",
    "320018208649"
)

for ($i = 0; $i -lt $newTracking.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 16)   # column P = ShipmentTracking
    # The tracking numbers are all-digit strings, so a plain .Value
    # assignment would be auto-parsed as a number by Excel. Force text
    # storage (matching the source data, which keeps these as strings),
    # then drop the style back to the sheet default so no stray
    # number-format is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $newTracking[$i]
    $cell.Style = "Normal"
}
